# Add a new order line item (row 9) to the order sheet.
#
# Every existing row stores all five columns (SKU, Name, Quantity, Cost Per,
# Total Cost) as plain text, even the numeric-looking ones - so the new row
# must match that: format the target cells as Text *before* assigning the
# values, otherwise Excel would auto-coerce the numeric-looking strings
# ("33576", "1", "69.75") into actual numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 9

$ws.Range("A9:E9").NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "33576"
$ws.Cells.Item($newRow, 2).Value = "Cup - Cold (12oz)"
$ws.Cells.Item($newRow, 3).Value = "1"
$ws.Cells.Item($newRow, 4).Value = "69.75"
$ws.Cells.Item($newRow, 5).Value = "69.75"
